# Insert a new weekly record for "Macroferia Regional de Talca - Zanahoria"
# at row 69, pushing the existing rows 69..160 down to 70..161.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(69).Insert()

$ws.Cells.Item(69, 1).Value = 5
$ws.Cells.Item(69, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(69, 3).Value = "Maule"
$ws.Cells.Item(69, 4).Value = 44413
$ws.Cells.Item(69, 5).Value = 7
$ws.Cells.Item(69, 6).Value = 100114013
$ws.Cells.Item(69, 7).Value = "Zanahoria"
$ws.Cells.Item(69, 8).Value = "Sin especificar"
$ws.Cells.Item(69, 9).Value = "Primera"
$ws.Cells.Item(69, 10).Value = 400
$ws.Cells.Item(69, 11).Value = 5000
$ws.Cells.Item(69, 12).Value = 5000
$ws.Cells.Item(69, 13).Value = 5000
$ws.Cells.Item(69, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(69, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(69, 16).Value = 250
$ws.Cells.Item(69, 17).Value = 20
$ws.Cells.Item(69, 18).Value = "Hortaliza"
